# Apply edits described by the commit:
# "adjusted rec selectivity and test regulation data"
#
# 1. Update fluke_min (col F) for rows 14-29, 41-55, 109-124 from 16.5 -> 17
# 2. Update fluke_min (col F) AND fluke_min_2019 (col R) for rows 56-64 from 16.5 -> 15
# 3. Update the sheet view / selection to reflect where the author was last working
#    (topLeftCell A103, active cell F109, selection F109:F124)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only column F (fluke_min) moves from 16.5 to 17
$rowsTo17 = @(14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,109,110,111,112,113,114,115,116,117,118,119,120,121,122,123,124)

foreach ($r in $rowsTo17) {
    $ws.Cells.Item($r, 6).Value = 17
}

# Rows where both column F (fluke_min) and column R (fluke_min_2019) move from 16.5 to 15
$rowsTo15 = @(56,57,58,59,60,61,62,63,64)

foreach ($r in $rowsTo15) {
    $ws.Cells.Item($r, 6).Value = 15
    $ws.Cells.Item($r, 18).Value = 15
}

# Update selection / scroll position to match author's last view
# (window scrolled so row 103 is the top visible row, F109:F124 selected)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 103
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F109:F124").Select()
